# [Word] (document) How to insert section breaks
# - Rename "word-manage-*" snippet ids to "word-document-manage-*"
# - Insert two new rows documenting Section.addNext / SectionCollection.addEven

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Rename the snippet id strings used across the table (column D) wherever they occur.
$rng = $ws.UsedRange
$rng.Replace("word-manage-comments", "word-document-manage-comments")
$rng.Replace("word-manage-footnotes", "word-document-manage-footnotes")
$rng.Replace("word-manage-change-tracking", "word-document-manage-change-tracking")

# 2) Insert a new row right before the "Section" / "getHeader" row (row 64) for the
#    new Section.addNext member.
$ws.Rows.Item(64).Insert()
$ws.Range("A64").Value2 = "Section"
$ws.Range("C64").Value2 = "class"
$ws.Range("D64").Value2 = "word-document-insert-section-breaks"
$ws.Range("E64").Value2 = "addNext"

# 3) Insert a new row right before the "Style" row (now row 67 after the previous
#    insert) for the new SectionCollection.addEven member.
$ws.Rows.Item(67).Insert()
$ws.Range("A67").Value2 = "SectionCollection"
$ws.Range("C67").Value2 = "class"
$ws.Range("D67").Value2 = "word-document-insert-section-breaks"
$ws.Range("E67").Value2 = "addEven"

# 4) Resize the table to include the two newly inserted rows.
$lo.Resize($ws.Range("A1:E69"))

# 5) Update the view state (frozen pane / selection) to match the authored file.
$ws.Range("D2").Select()
$actWin = $excel.ActiveWindow
$actWin.ScrollRow = 1
$actWin.FreezePanes = $false
$ws.Range("A2").Select()
$actWin.FreezePanes = $true
$ws.Range("D2").Select()
